$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# Update data values (counts) for rows 2-5
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 3

$ws.Range("B3").Value = 3
$ws.Range("D3").Value = 9

$ws.Range("B4").Value = 8
$ws.Range("D4").Value = 35

$ws.Range("B5").Value = 6
$ws.Range("D5").Value = 27

# Update computed percentages
$ws.Range("C2").Value = 5.555555555555555
$ws.Range("E2").Value = 4.054054054054054

$ws.Range("C3").Value = 16.66666666666666
$ws.Range("E3").Value = 12.16216216216216

$ws.Range("C4").Value = 44.44444444444444
$ws.Range("E4").Value = 47.2972972972973

$ws.Range("C5").Value = 33.33333333333333
$ws.Range("E5").Value = 36.48648648648648
